$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'30.661.47"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.54%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'2.116.27"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  +0.20%  "
$c.Style = "Normal"

$c = $ws.Range("D4")
$c.Value = "'1.013"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  +1.01%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'339.27"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  +1.36%  "
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'  +1.02%  "
$c.Style = "Normal"

$c = $ws.Range("D7")
$c.Value = "'0.5256"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +0.08%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'0.4513"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +0.03%  "
$c.Style = "Normal"

$c = $ws.Range("D9")
$c.Value = "'53.90"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  +0.70%  "
$c.Style = "Normal"

$c = $ws.Range("D10")
$c.Value = "'0.09098"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +1.04%  "
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'  +0.11%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'24.38"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -0.59%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'2.121.54"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +0.42%  "
$c.Style = "Normal"

$c = $ws.Range("D14")
$c.Value = "'6.808"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +0.23%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'8.082"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  +3.22%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'97.78"
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +0.98%  "
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'  +3.11%  "
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'  +1.05%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'0.06702"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  +1.16%  "
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'  -0.13%  "
$c.Style = "Normal"

$c = $ws.Range("D21")
$c.Value = "'1.012"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +1.01%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'6.413"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  +1.56%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'30.748.35"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  +0.67%  "
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'  +3.68%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'2.383"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  +1.11%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'2.368.90"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  +0.43%  "
$c.Style = "Normal"

$c = $ws.Range("D27")
$c.Value = "'22.38"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -0.20%  "
$c.Style = "Normal"

$c = $ws.Range("D28")
$c.Value = "'164.93"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  +0.91%  "
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "'  -1.52%  "
$c.Style = "Normal"

$c = $ws.Range("D30")
$c.Value = "'134.99"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +1.63%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'1.200"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -0.29%  "
$c.Style = "Normal"

$c = $ws.Range("D32")
$c.Value = "'0.1075"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  -0.01%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'6.378"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  +3.32%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'1.637"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  -1.57%  "
$c.Style = "Normal"

$c = $ws.Range("D35")
$c.Value = "'3.946"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  +0.22%  "
$c.Style = "Normal"

$c = $ws.Range("D36")
$c.Value = "'10.37"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  -2.92%  "
$c.Style = "Normal"

$c = $ws.Range("D37")
$c.Value = "'5.917"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  +6.46%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'0.02656"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +2.91%  "
$c.Style = "Normal"

$c = $ws.Range("D39")
$c.Value = "'0.06835"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -0.12%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.2323"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  +0.99%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'12.60"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -1.89%  "
$c.Style = "Normal"

$c = $ws.Range("D42")
$c.Value = "'0.6886"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -0.95%  "
$c.Style = "Normal"

$c = $ws.Range("D43")
$c.Value = "'1.263"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +1.07%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'14.97"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +6.53%  "
$c.Style = "Normal"

$c = $ws.Range("D45")
$c.Value = "'0.6444"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  +0.40%  "
$c.Style = "Normal"

$c = $ws.Range("D46")
$c.Value = "'2.323"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -3.63%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'0.00000000371"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  +14.88%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'3.707"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  +1.10%  "
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'  +0.27%  "
$c.Style = "Normal"

$c = $ws.Range("D50")
$c.Value = "'0.07314"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +3.17%  "
$c.Style = "Normal"

$c = $ws.Range("D51")
$c.Value = "'82.83"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -0.74%  "
$c.Style = "Normal"
